# Apply the commit's data edits to the proteomics worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the "significant" boolean flags (column B) for the rows that
# changed between the previous and the new analysis run.
$ws.Range("B5:B9").Value = $false
$ws.Range("B13:B19").Value = $false
$ws.Range("B22:B28").Value = $false
$ws.Range("B31:B37").Value = $false
$ws.Range("B41:B44").Value = $false
$ws.Range("B49:B53").Value = $false
$ws.Range("B57").Value = $true
$ws.Range("B66").Value = $true
$ws.Range("B73").Value = $false
$ws.Range("B75").Value = $true
$ws.Range("B82").Value = $false
$ws.Range("B86:B90").Value = $false
$ws.Range("B95:B96").Value = $false
$ws.Range("B102:B103").Value = $true
$ws.Range("B111:B112").Value = $true
$ws.Range("B120:B121").Value = $true
$ws.Range("B133:B134").Value = $false
$ws.Range("B142").Value = $false
$ws.Range("B150:B153").Value = $false
$ws.Range("B159:B162").Value = $false
$ws.Range("B169:B171").Value = $false

# Re-select the data range so the workbook records the new active
# cell / selection state on the sheet view.
$ws.Range("A2:I181").Select()
